$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (3 and 4) for "Arveja Verde" had their weekly report
# values swapped: row 3 now holds the later week's figures and row 4 the
# earlier week's figures.

$row3 = @{
    D = 44827
    J = 300
    K = 30000
    L = 31000
    M = 30500
    P = 1220
}

$row4 = @{
    D = 44414
    J = 500
    K = 31000
    L = 32000
    M = 31500
    P = 1260
}

foreach ($col in $row3.Keys) {
    $ws.Range("$col`3").Value = $row3[$col]
}

foreach ($col in $row4.Keys) {
    $ws.Range("$col`4").Value = $row4[$col]
}
